$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45998
$ws.Range("B2").Value = 33.32
$ws.Range("C2").Value = 24
$ws.Range("D2").Value = 22.71
$ws.Range("E2").Value = 17.68
$ws.Range("F2").Value = 11
$ws.Range("G2").Value = 14.53
$ws.Range("H2").Value = 18.37
$ws.Range("I2").Value = 18.99
$ws.Range("J2").Value = 19.92
$ws.Range("K2").Value = 16.33
$ws.Range("L2").Value = 11.14
$ws.Range("M2").Value = 6.27
$ws.Range("N2").Value = 0.74
$ws.Range("O2").Value = 0.5
$ws.Range("P2").Value = 0.8100000000000001
$ws.Range("Q2").Value = 4.55
$ws.Range("R2").Value = 15.08
$ws.Range("S2").Value = 49.92
$ws.Range("T2").Value = 73.63
$ws.Range("U2").Value = 83.03
$ws.Range("V2").Value = 89.59
$ws.Range("W2").Value = 91.14
$ws.Range("X2").Value = 80.3
$ws.Range("Y2").Value = 71.59999999999999
$ws.Range("Z2").Value = 32.3
$ws.Range("AB2").Value = 83.16
$ws.Range("AD2").Value = 90.36
$ws.Range("AF2").Value = 78.33
